$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate/shift existing rows per diff ---
$ws.Cells.Item(23,6).Value = "Dravinja"
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = "Ilirija"
$ws.Cells.Item(23,9).Value = 0
$ws.Cells.Item(23,10).Value = 3.12
$ws.Cells.Item(23,11).Value = "20/08/2023 08:59"
$ws.Cells.Item(23,12).Value = 3.08
$ws.Cells.Item(23,13).Value = "20/08/2023 11:20"
$ws.Cells.Item(23,14).Value = 3.27
$ws.Cells.Item(23,15).Value = "20/08/2023 08:59"
$ws.Cells.Item(23,16).Value = 3.4
$ws.Cells.Item(23,17).Value = "20/08/2023 15:35"
$ws.Cells.Item(23,18).Value = 2.15
$ws.Cells.Item(23,19).Value = "20/08/2023 08:59"
$ws.Cells.Item(23,20).Value = 2.12
$ws.Cells.Item(23,21).Value = "20/08/2023 14:26"
$ws.Cells.Item(23,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/dravinja-ilirija/zRWDoyR5/"

$ws.Cells.Item(24,6).Value = "Bilje"
$ws.Cells.Item(24,7).Value = 1
$ws.Cells.Item(24,8).Value = "Tabor Sezana"
$ws.Cells.Item(24,9).Value = 1
$ws.Cells.Item(24,10).Value = 1.52
$ws.Cells.Item(24,11).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,12).Value = 1.47
$ws.Cells.Item(24,13).Value = "20/08/2023 17:19"
$ws.Cells.Item(24,14).Value = 4.26
$ws.Cells.Item(24,15).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,16).Value = 5.59
$ws.Cells.Item(24,17).Value = "20/08/2023 17:19"
$ws.Cells.Item(24,18).Value = 5.03
$ws.Cells.Item(24,19).Value = "20/08/2023 09:00"
$ws.Cells.Item(24,20).Value = 4.39
$ws.Cells.Item(24,21).Value = "20/08/2023 17:25"
$ws.Cells.Item(24,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-tabor-sezana/6ouckJRt/"

$ws.Cells.Item(25,6).Value = "NK Krka"
$ws.Cells.Item(25,7).Value = 1
$ws.Cells.Item(25,8).Value = "NK Bistrica"
$ws.Cells.Item(25,9).Value = 2
$ws.Cells.Item(25,10).Value = 1.88
$ws.Cells.Item(25,11).Value = "19/08/2023 05:42"
$ws.Cells.Item(25,12).Value = 2.26
$ws.Cells.Item(25,13).Value = "20/08/2023 17:19"
$ws.Cells.Item(25,14).Value = 3.4
$ws.Cells.Item(25,15).Value = "19/08/2023 05:42"
$ws.Cells.Item(25,16).Value = 3.33
$ws.Cells.Item(25,17).Value = "20/08/2023 17:19"
$ws.Cells.Item(25,18).Value = 3.31
$ws.Cells.Item(25,19).Value = "19/08/2023 05:42"
$ws.Cells.Item(25,20).Value = 2.93
$ws.Cells.Item(25,21).Value = "20/08/2023 17:19"
$ws.Cells.Item(25,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nk-krka-bistrica/n1zIpetC/"

$ws.Cells.Item(29,6).Value = "Fuzinar"
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = "Grosuplje"
$ws.Cells.Item(29,9).Value = 4
$ws.Cells.Item(29,10).Value = 3.24
$ws.Cells.Item(29,11).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,12).Value = 4.09
$ws.Cells.Item(29,13).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,14).Value = 3.15
$ws.Cells.Item(29,15).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,16).Value = 3.69
$ws.Cells.Item(29,17).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,18).Value = 2.03
$ws.Cells.Item(29,19).Value = "04/08/2023 05:42"
$ws.Cells.Item(29,20).Value = 1.75
$ws.Cells.Item(29,21).Value = "30/08/2023 15:38"
$ws.Cells.Item(29,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-grosuplje/OMD8pzPE/"

$ws.Cells.Item(30,6).Value = "Ilirija"
$ws.Cells.Item(30,7).Value = 0
$ws.Cells.Item(30,8).Value = "Rudar"
$ws.Cells.Item(30,9).Value = 2
$ws.Cells.Item(30,10).Value = 1.81
$ws.Cells.Item(30,11).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,12).Value = 2.13
$ws.Cells.Item(30,13).Value = "30/08/2023 16:51"
$ws.Cells.Item(30,14).Value = 3.51
$ws.Cells.Item(30,15).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,16).Value = 3.62
$ws.Cells.Item(30,17).Value = "30/08/2023 16:21"
$ws.Cells.Item(30,18).Value = 3.42
$ws.Cells.Item(30,19).Value = "05/08/2023 05:42"
$ws.Cells.Item(30,20).Value = 2.95
$ws.Cells.Item(30,21).Value = "30/08/2023 16:51"
$ws.Cells.Item(30,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/ilirija-rudar/pv9TMWW7/"

$ws.Cells.Item(33,6).Value = "Primorje"
$ws.Cells.Item(33,7).Value = 3
$ws.Cells.Item(33,8).Value = "Nafta"
$ws.Cells.Item(33,9).Value = 0
$ws.Cells.Item(33,10).Value = 2.32
$ws.Cells.Item(33,11).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,12).Value = 2.33
$ws.Cells.Item(33,13).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,14).Value = 3.25
$ws.Cells.Item(33,15).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,16).Value = 3.38
$ws.Cells.Item(33,17).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,18).Value = 2.63
$ws.Cells.Item(33,19).Value = "01/09/2023 04:43"
$ws.Cells.Item(33,20).Value = 2.77
$ws.Cells.Item(33,21).Value = "02/09/2023 16:27"
$ws.Cells.Item(33,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/primorje-nafta/GKUtkAXc/"

$ws.Cells.Item(34,6).Value = "Tabor Sezana"
$ws.Cells.Item(34,7).Value = 1
$ws.Cells.Item(34,8).Value = "Rudar"
$ws.Cells.Item(34,9).Value = 1
$ws.Cells.Item(34,10).Value = 4.32
$ws.Cells.Item(34,11).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,12).Value = 3.42
$ws.Cells.Item(34,13).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,14).Value = 3.88
$ws.Cells.Item(34,15).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,16).Value = 3.99
$ws.Cells.Item(34,17).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,18).Value = 1.66
$ws.Cells.Item(34,19).Value = "02/09/2023 14:10"
$ws.Cells.Item(34,20).Value = 1.85
$ws.Cells.Item(34,21).Value = "02/09/2023 16:28"
$ws.Cells.Item(34,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-rudar/UNYplUn4/"

$ws.Cells.Item(35,6).Value = "Fuzinar"
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = "Jadran Dekani"
$ws.Cells.Item(35,9).Value = 3
$ws.Cells.Item(35,10).Value = 2.35
$ws.Cells.Item(35,11).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,12).Value = 2.84
$ws.Cells.Item(35,13).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,14).Value = 3.12
$ws.Cells.Item(35,15).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,16).Value = 3.36
$ws.Cells.Item(35,17).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,18).Value = 2.68
$ws.Cells.Item(35,19).Value = "01/09/2023 04:43"
$ws.Cells.Item(35,20).Value = 2.3
$ws.Cells.Item(35,21).Value = "02/09/2023 16:21"
$ws.Cells.Item(35,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-jadran-dekani/YeuxjjIi/"

$ws.Cells.Item(36,6).Value = "ND Gorica"
$ws.Cells.Item(36,7).Value = 3
$ws.Cells.Item(36,8).Value = "Dravinja"
$ws.Cells.Item(36,9).Value = 0
$ws.Cells.Item(36,10).Value = 1.58
$ws.Cells.Item(36,11).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,12).Value = 1.44
$ws.Cells.Item(36,13).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,14).Value = 3.99
$ws.Cells.Item(36,15).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,16).Value = 4.38
$ws.Cells.Item(36,17).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,18).Value = 4.75
$ws.Cells.Item(36,19).Value = "02/09/2023 14:10"
$ws.Cells.Item(36,20).Value = 6.18
$ws.Cells.Item(36,21).Value = "02/09/2023 14:41"
$ws.Cells.Item(36,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nd-gorica-dravinja/f5tYjW2o/"

$ws.Cells.Item(44,6).Value = "Nafta"
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = "Tabor Sezana"
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 1.19
$ws.Cells.Item(44,11).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,12).Value = 1.29
$ws.Cells.Item(44,13).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,14).Value = 6.85
$ws.Cells.Item(44,15).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,16).Value = 5.9
$ws.Cells.Item(44,17).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,18).Value = 8.87
$ws.Cells.Item(44,19).Value = "09/09/2023 13:42"
$ws.Cells.Item(44,20).Value = 7.25
$ws.Cells.Item(44,21).Value = "09/09/2023 16:23"
$ws.Cells.Item(44,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/nafta-tabor-sezana/hMwQvQ9j/"

$ws.Cells.Item(46,6).Value = "Rudar"
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = "Beltinci"
$ws.Cells.Item(46,9).Value = 1
$ws.Cells.Item(46,10).Value = 2.41
$ws.Cells.Item(46,11).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,12).Value = 3.4
$ws.Cells.Item(46,13).Value = "09/09/2023 16:15"
$ws.Cells.Item(46,14).Value = 3.23
$ws.Cells.Item(46,15).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,16).Value = 3.56
$ws.Cells.Item(46,17).Value = "09/09/2023 16:14"
$ws.Cells.Item(46,18).Value = 2.54
$ws.Cells.Item(46,19).Value = "08/09/2023 04:42"
$ws.Cells.Item(46,20).Value = 1.91
$ws.Cells.Item(46,21).Value = "09/09/2023 16:15"
$ws.Cells.Item(46,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/rudar-beltinci/vJZLu6fp/"

$ws.Cells.Item(47,6).Value = "Bilje"
$ws.Cells.Item(47,7).Value = 2
$ws.Cells.Item(47,8).Value = "Tolmin"
$ws.Cells.Item(47,9).Value = 1
$ws.Cells.Item(47,10).Value = 1.54
$ws.Cells.Item(47,11).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,12).Value = 1.49
$ws.Cells.Item(47,13).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,14).Value = 3.85
$ws.Cells.Item(47,15).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,16).Value = 4.38
$ws.Cells.Item(47,17).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,18).Value = 4.49
$ws.Cells.Item(47,19).Value = "08/09/2023 04:42"
$ws.Cells.Item(47,20).Value = 5.37
$ws.Cells.Item(47,21).Value = "09/09/2023 16:22"
$ws.Cells.Item(47,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/bilje-tolmin/MwmfpnnT/"


# --- Append new rows 66-71 ---
$ws.Range("A65:V65").Copy()
$ws.Range("A66:V66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(66,1).Value = 65
$ws.Cells.Item(66,2).Value = "slovenia"
$ws.Cells.Item(66,3).Value = "2-snl"
$ws.Cells.Item(66,4).Value = "2023-2024"
$ws.Cells.Item(66,5).Value = 45192.66666666666
$ws.Cells.Item(66,6).Value = "Fuzinar"
$ws.Cells.Item(66,7).Value = 3
$ws.Cells.Item(66,8).Value = "Ilirija"
$ws.Cells.Item(66,9).Value = 2
$ws.Cells.Item(66,10).Value = 2.44
$ws.Cells.Item(66,11).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,12).Value = 2.45
$ws.Cells.Item(66,13).Value = "23/09/2023 15:58"
$ws.Cells.Item(66,14).Value = 3.32
$ws.Cells.Item(66,15).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,16).Value = 3.7
$ws.Cells.Item(66,17).Value = "23/09/2023 15:59"
$ws.Cells.Item(66,18).Value = 2.45
$ws.Cells.Item(66,19).Value = "22/09/2023 03:13"
$ws.Cells.Item(66,20).Value = 2.46
$ws.Cells.Item(66,21).Value = "23/09/2023 15:58"
$ws.Cells.Item(66,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/fuzinar-ilirija/dY8ySqV8/"

$ws.Range("A66:V66").Copy()
$ws.Range("A67:V67").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(67,1).Value = 66
$ws.Cells.Item(67,2).Value = "slovenia"
$ws.Cells.Item(67,3).Value = "2-snl"
$ws.Cells.Item(67,4).Value = "2023-2024"
$ws.Cells.Item(67,5).Value = 45192.66666666666
$ws.Cells.Item(67,6).Value = "Primorje"
$ws.Cells.Item(67,7).Value = 2
$ws.Cells.Item(67,8).Value = "NK Bistrica"
$ws.Cells.Item(67,9).Value = 2
$ws.Cells.Item(67,10).Value = 1.84
$ws.Cells.Item(67,11).Value = "22/09/2023 03:13"
$ws.Cells.Item(67,12).Value = 1.75
$ws.Cells.Item(67,13).Value = "23/09/2023 15:46"
$ws.Cells.Item(67,14).Value = 3.4
$ws.Cells.Item(67,15).Value = "22/09/2023 03:13"
$ws.Cells.Item(67,16).Value = 3.65
$ws.Cells.Item(67,17).Value = "23/09/2023 15:46"
$ws.Cells.Item(67,18).Value = 3.42
$ws.Cells.Item(67,19).Value = "22/09/2023 03:13"
$ws.Cells.Item(67,20).Value = 4.15
$ws.Cells.Item(67,21).Value = "23/09/2023 15:46"
$ws.Cells.Item(67,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/primorje-bistrica/je8uR3pF/"

$ws.Range("A67:V67").Copy()
$ws.Range("A68:V68").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(68,1).Value = 67
$ws.Cells.Item(68,2).Value = "slovenia"
$ws.Cells.Item(68,3).Value = "2-snl"
$ws.Cells.Item(68,4).Value = "2023-2024"
$ws.Cells.Item(68,5).Value = 45192.66666666666
$ws.Cells.Item(68,6).Value = "Tabor Sezana"
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = "NK Krka"
$ws.Cells.Item(68,9).Value = 2
$ws.Cells.Item(68,10).Value = 3.75
$ws.Cells.Item(68,11).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,12).Value = 4.99
$ws.Cells.Item(68,13).Value = "23/09/2023 15:47"
$ws.Cells.Item(68,14).Value = 3.56
$ws.Cells.Item(68,15).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,16).Value = 4.23
$ws.Cells.Item(68,17).Value = "23/09/2023 15:47"
$ws.Cells.Item(68,18).Value = 1.72
$ws.Cells.Item(68,19).Value = "22/09/2023 03:13"
$ws.Cells.Item(68,20).Value = 1.54
$ws.Cells.Item(68,21).Value = "23/09/2023 15:47"
$ws.Cells.Item(68,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tabor-sezana-nk-krka/pf1lPsFR/"

$ws.Range("A68:V68").Copy()
$ws.Range("A69:V69").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(69,1).Value = 68
$ws.Cells.Item(69,2).Value = "slovenia"
$ws.Cells.Item(69,3).Value = "2-snl"
$ws.Cells.Item(69,4).Value = "2023-2024"
$ws.Cells.Item(69,5).Value = 45193.66666666666
$ws.Cells.Item(69,6).Value = "Beltinci"
$ws.Cells.Item(69,7).Value = 3
$ws.Cells.Item(69,8).Value = "Dravinja"
$ws.Cells.Item(69,9).Value = 0
$ws.Cells.Item(69,10).Value = 1.4
$ws.Cells.Item(69,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,12).Value = 1.38
$ws.Cells.Item(69,13).Value = "24/09/2023 15:43"
$ws.Cells.Item(69,14).Value = 4.36
$ws.Cells.Item(69,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,16).Value = 4.76
$ws.Cells.Item(69,17).Value = "24/09/2023 15:54"
$ws.Cells.Item(69,18).Value = 5.39
$ws.Cells.Item(69,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(69,20).Value = 6.65
$ws.Cells.Item(69,21).Value = "24/09/2023 15:54"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/beltinci-dravinja/n3mGK1Nr/"

$ws.Range("A69:V69").Copy()
$ws.Range("A70:V70").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(70,1).Value = 69
$ws.Cells.Item(70,2).Value = "slovenia"
$ws.Cells.Item(70,3).Value = "2-snl"
$ws.Cells.Item(70,4).Value = "2023-2024"
$ws.Cells.Item(70,5).Value = 45193.66666666666
$ws.Cells.Item(70,6).Value = "Tolmin"
$ws.Cells.Item(70,7).Value = 2
$ws.Cells.Item(70,8).Value = "Jadran Dekani"
$ws.Cells.Item(70,9).Value = 2
$ws.Cells.Item(70,10).Value = 3.15
$ws.Cells.Item(70,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,12).Value = 3.52
$ws.Cells.Item(70,13).Value = "24/09/2023 15:42"
$ws.Cells.Item(70,14).Value = 3.16
$ws.Cells.Item(70,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,16).Value = 3.17
$ws.Cells.Item(70,17).Value = "24/09/2023 15:41"
$ws.Cells.Item(70,18).Value = 2.06
$ws.Cells.Item(70,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(70,20).Value = 2.06
$ws.Cells.Item(70,21).Value = "24/09/2023 15:41"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/tolmin-jadran-dekani/GxhKJLxk/"

$ws.Range("A70:V70").Copy()
$ws.Range("A71:V71").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(71,1).Value = 70
$ws.Cells.Item(71,2).Value = "slovenia"
$ws.Cells.Item(71,3).Value = "2-snl"
$ws.Cells.Item(71,4).Value = "2023-2024"
$ws.Cells.Item(71,5).Value = 45193.66666666666
$ws.Cells.Item(71,6).Value = "Triglav"
$ws.Cells.Item(71,7).Value = 0
$ws.Cells.Item(71,8).Value = "ND Gorica"
$ws.Cells.Item(71,9).Value = 3
$ws.Cells.Item(71,10).Value = 3.01
$ws.Cells.Item(71,11).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,12).Value = 3.77
$ws.Cells.Item(71,13).Value = "24/09/2023 15:32"
$ws.Cells.Item(71,14).Value = 3.29
$ws.Cells.Item(71,15).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,16).Value = 3.27
$ws.Cells.Item(71,17).Value = "24/09/2023 15:39"
$ws.Cells.Item(71,18).Value = 2.03
$ws.Cells.Item(71,19).Value = "23/09/2023 03:12"
$ws.Cells.Item(71,20).Value = 1.94
$ws.Cells.Item(71,21).Value = "24/09/2023 15:39"
$ws.Cells.Item(71,22).Value = "https://www.betexplorer.com/football/slovenia/2-snl/triglav-nd-gorica/IsAXSPF2/"

